$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1599
$ws1.Range("F4").Value = 844
$ws1.Range("F5").Value = 241
$ws1.Range("F6").Value = 72
$ws1.Range("F7").Value = 1153
$ws1.Range("F8").Value = 757
$ws1.Range("F9").Value = 800
$ws1.Range("F10").Value = 1456
$ws1.Range("F11").Value = 290
$ws1.Range("F12").Value = 1037
$ws1.Range("F14").Value = 67
$ws1.Range("F17").Value = 481
$ws1.Range("F18").Value = 33
$ws1.Range("F19").Value = 32
$ws1.Range("F22").Value = 299
$ws1.Range("F23").Value = 556
$ws1.Range("F24").Value = 571
$ws1.Range("F26").Value = 248
$ws1.Range("F27").Value = 184
$ws1.Range("F28").Value = 372
$ws1.Range("D23").Value = "广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心"

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 1006
$ws2.Range("F5").Value = 272
$ws2.Range("F9").Value = 590
$ws2.Range("F10").Value = 85

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 248
$ws4.Range("F3").Value = 248
$ws4.Range("F4").Value = 1599
$ws4.Range("F6").Value = 844
$ws4.Range("F7").Value = 241
$ws4.Range("F8").Value = 1006
$ws4.Range("F9").Value = 72
$ws4.Range("F10").Value = 1153
$ws4.Range("F11").Value = 757
$ws4.Range("F12").Value = 800
$ws4.Range("F13").Value = 1456
$ws4.Range("F14").Value = 290
$ws4.Range("F15").Value = 1037
$ws4.Range("F17").Value = 67
$ws4.Range("F20").Value = 481
$ws4.Range("F21").Value = 33
$ws4.Range("F22").Value = 32
$ws4.Range("F25").Value = 272
$ws4.Range("F27").Value = 299
$ws4.Range("F31").Value = 556
$ws4.Range("F32").Value = 571
$ws4.Range("F34").Value = 248
$ws4.Range("F36").Value = 184
$ws4.Range("F37").Value = 590
$ws4.Range("F38").Value = 85
$ws4.Range("F39").Value = 85
$ws4.Range("F41").Value = 372
$ws4.Range("D31").Value = "广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心"
